$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new row at position 7 ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = 11
$ws.Range("B7").Value = "Vega Monumental Concepción"
$ws.Range("C7").Value = "Bíobío"
$ws.Range("D7").Value = 44425
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 100112012
$ws.Range("G7").Value = "Espinaca"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 60
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 7000
$ws.Range("M7").Value = 6750
$ws.Range("N7").Value = "`$/cuna 10 kilos"
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 675
$ws.Range("Q7").Value = 10
$ws.Range("R7").Value = "Hortaliza"

# --- Insert new row at position 16 (post first insert) ---
$ws.Rows.Item(16).Insert()
$ws.Range("A16").Value = 11
$ws.Range("B16").Value = "Vega Monumental Concepción"
$ws.Range("C16").Value = "Bíobío"
$ws.Range("D16").Value = 44421
$ws.Range("E16").Value = 8
$ws.Range("F16").Value = 100112012
$ws.Range("G16").Value = "Espinaca"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 7000
$ws.Range("L16").Value = 7500
$ws.Range("M16").Value = 7250
$ws.Range("N16").Value = "`$/cuna 10 kilos"
$ws.Range("O16").Value = "Región Metropolitana"
$ws.Range("P16").Value = 725
$ws.Range("Q16").Value = 10
$ws.Range("R16").Value = "Hortaliza"

# Ensure the date columns keep the date number format used elsewhere in column D
$ws.Range("D7").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Edit complete"
